# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-holding detail, same layout as
#    the other quarterly sheets) right before the "总计" (summary) sheet.
# 2) Prepend a new "2022-Q1" row to the "总计" summary sheet and renumber
#    the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
# Copy the "2021-Q4" sheet (same column layout/styles) so the new sheet
# inherits identical formatting, then drop it in right before "总计".
$template = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$totalIdx = $total.Index

$template.Copy($total)

$newSheet = $wb.Worksheets.Item($totalIdx)
$newSheet.Name = "2022-Q1"

# Fund holdings for 2022-Q1 (index, code, name, fund size, total stock
# position, position ratio, held market value (100M yuan), position rank).
$data = @(
    @(0, "009548", "汇添富中盘价值精选混合A",          "161.71", "84.81", "4.06", "6.5654", 7),
    @(1, "900010", "中信卓越成长两年持有期混合A",      "133.02", "93.07", "3.60", "4.7887", 5),
    @(2, "900090", "中信卓越成长两年持有期混合B",      "86.95",  "93.07", "3.60", "3.1302", 5),
    @(3, "009549", "汇添富中盘价值精选混合C",          "21.00",  "84.81", "4.06", "0.8526", 7),
    @(4, "001371", "富国沪港深价值精选灵活配置混合A",  "40.15",  "68.33", "1.97", "0.7910", 6),
    @(5, "900100", "中信卓越成长两年持有期混合C",      "6.91",   "93.07", "3.60", "0.2488", 5),
    @(6, "011131", "富国沪港深价值精选灵活配置混合C",  "0.42",   "68.33", "1.97", "0.0083", 6)
)

$lastTemplateRow = 7
$r = 2
foreach ($row in $data) {
    if ($r -gt $lastTemplateRow) {
        # Row 8 doesn't exist in the 6-row template yet: clone the index
        # column's style from the row above before filling values in.
        $newSheet.Range("A" + ($r - 1)).Copy()
        $newSheet.Range("A" + $r).PasteSpecial(-4122)
    }

    $newSheet.Range("A$r").Value2 = $row[0]

    # Force text storage (matches the source sheets, which store these
    # numeric-looking columns as text) via a leading apostrophe, then
    # strip the resulting "quote prefix" style back to the default so
    # the cells don't pick up a new/different style index.
    $newSheet.Range("B$r").Value = "'" + $row[1]
    $newSheet.Range("C$r").Value = "'" + $row[2]
    $newSheet.Range("D$r").Value = "'" + $row[3]
    $newSheet.Range("E$r").Value = "'" + $row[4]
    $newSheet.Range("F$r").Value = "'" + $row[5]
    $newSheet.Range("G$r").Value = "'" + $row[6]
    $newSheet.Range("B$r`:G$r").ClearFormats()

    $newSheet.Range("H$r").Value2 = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Part 2: prepend a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

# The inserted row doesn't inherit the index column's border/alignment
# style automatically - copy it from the row beneath (old row 2).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value2 = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value2 = 7
$totalSheet.Range("D2").Value2 = 16.38
$totalSheet.Range("B2:D2").ClearFormats()
# ClearFormats also resets B2/C2/D2 to default style - fine, they match
# the rest of the data rows (no explicit style).

# Renumber the index column (A) for the rows that got shifted down.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value2 = $r - 2
}
